$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "tablet"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "gear"

$ws.Range("B12").Select()
